# Update the "NIEM Xpath" mapping column (C2:C14): the service/document
# root name changed from evalsr-doc:PersonEvaluationSearchRequest to
# phisr-doc:PersonHealthInformationSearchRequest (and the companion
# extension namespace prefix evalsr-ext -> phisr-ext).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonGivenName"
$ws.Range("C3").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonMiddleName"
$ws.Range("C4").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonName/nc30:PersonSurName"
$ws.Range("C5").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/nc30:PersonBirthDate/nc30:Date"
$ws.Range("C6").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonRaceCode"
$ws.Range("C7").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonEthnicityCode"
$ws.Range("C8").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/jxdm51:PersonSexCode"
$ws.Range("C9").Value  = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/phisr-ext:PersonTemporaryIdentification/nc30:IdentificationID"
$ws.Range("C10").Value = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonGivenName"
$ws.Range("C11").Value = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonMiddleName"
$ws.Range("C12").Value = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonName/nc30:PersonSurName"
$ws.Range("C13").Value = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/nc30:PersonBirthDate/nc30:Date"
$ws.Range("C14").Value = "/phisr-doc:PersonHealthInformationSearchRequest/nc30:Identity[@structures:id=/phisr-doc:PersonHealthInformationSearchRequest/nc30:PersonAliasIdentityAssociation[nc30:Person/@structures:ref=/phisr-doc:PersonHealthInformationSearchRequest/nc30:Person/@structures:id]/nc30:Identity/@structures:ref]/nc30:IdentityPersonRepresentation/jxdm51:PersonSexCode"

# The sheet's last recorded selection moves from C2 to C16 (below the data).
$ws.Range("C16").Select()
